$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. New "Average" column (J) on the left-hand (Mean) table -- this must be
#    written first so that the shared-string table gets "Average" as the
#    first newly-added unique string (matching the original author's edit
#    order: Average, Volunteer, MAX:).
# ---------------------------------------------------------------------------
$ws.Range("J3").Value = "Average"
$ws.Range("J3").Font.Bold = $true

$ws.Range("J4").Formula = "=AVERAGE(B4:I4)"
$ws.Range("J5:J20").Formula = "=AVERAGE(B5:I5)"

# "Average" column (T) on the right-hand (StD) table -- reuses the shared
# string created above, so it does not add a new entry.
$ws.Range("T3").Value = "Average"

$ws.Range("T4").Formula = "=AVERAGE(L4:S4)"
$ws.Range("T5:T20").Formula = "=AVERAGE(L5:S5)"

# ---------------------------------------------------------------------------
# 2. New "Volunteer" label above column A (row 3)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "Volunteer"

# ---------------------------------------------------------------------------
# 3. "MAX:" summary row (row 21)
# ---------------------------------------------------------------------------
$ws.Range("A21").Value = "MAX:"

$ws.Range("B21").Formula = "=MAX(B4:B20)"
$ws.Range("C21:I21").Formula = "=MAX(C4:C20)"
$ws.Range("J21").Formula = "=MAX(J4:J20)"

# ---------------------------------------------------------------------------
# 4. View: selection / scroll position
# ---------------------------------------------------------------------------
$ws.Range("A21:B21").Select()
$excel.ActiveWindow.ScrollRow = 3

# ---------------------------------------------------------------------------
# 5. Page setup (portrait) + charts / drawing
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1

# Chart 1: clustered column chart, one series per electrode (B..J)
$chartObj1 = $ws.ChartObjects().Add(114299, 5981699, 6080125, 3667125)
$chartObj1.Name = "Chart 3"
$chart1 = $chartObj1.Chart
$chart1.ChartType = 51
$series1 = $chart1.SeriesCollection()
$cols = @("B", "C", "D", "E", "F", "G", "H", "I", "J")
foreach ($col in $cols) {
    $ser = $series1.NewSeries()
    $ser.Name = "=Sheet1!`$$col`$3"
    $ser.Values = $ws.Range($col + "4:" + $col + "20")
}
$chart1.HasLegend = $true
$chart1.Legend.Position = -4152
$chart1.HasTitle = $false

# Chart 2: clustered column chart, single "Average" series
$chartObj2 = $ws.ChartObjects().Add(1162050, 3295650, 2486025, 1400175)
$chartObj2.Name = "Chart 6"
$chart2 = $chartObj2.Chart
$chart2.ChartType = 51
$ser2 = $chart2.SeriesCollection().NewSeries()
$ser2.Name = "=Sheet1!`$J`$3"
$ser2.Values = $ws.Range("J4:J20")
$chart2.HasLegend = $true
$chart2.Legend.Position = -4152
$chart2.HasTitle = $true
$chart2.ChartTitle.Text = ""
